$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.831.68'
$ws.Cells.Item(2, 5).Value = '  -2.27%  '
$ws.Cells.Item(3, 4).Value = '3.832.74'
$ws.Cells.Item(4, 5).Value = '  -0.03%  '
$ws.Cells.Item(5, 4).Value = "'599.43"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -2.58%  '
$ws.Cells.Item(6, 4).Value = "'178.71"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.66%  '
$ws.Cells.Item(7, 4).Value = '3.829.15'
$ws.Cells.Item(7, 5).Value = '  +1.54%  '
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 4).Value = "'0.529"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +0.36%  '
$ws.Cells.Item(10, 4).Value = "'0.162"
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -3.61%  '
$ws.Cells.Item(11, 4).Value = "'6.20"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -5.44%  '
$ws.Cells.Item(12, 4).Value = "'0.467"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -3.36%  '
$ws.Cells.Item(13, 4).Value = "'38.73"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -3.81%  '
$ws.Cells.Item(14, 5).Value = '  -3.53%  '
$ws.Cells.Item(15, 4).Value = '4.479.14'
$ws.Cells.Item(15, 5).Value = '  +1.76%  '
$ws.Cells.Item(16, 4).Value = '3.834.92'
$ws.Cells.Item(16, 5).Value = '  +1.77%  '
$ws.Cells.Item(17, 4).Value = '68.005.65'
$ws.Cells.Item(17, 5).Value = '  -2.26%  '
$ws.Cells.Item(18, 4).Value = "'7.24"
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -3.25%  '
$ws.Cells.Item(19, 5).Value = '  -4.53%  '
$ws.Cells.Item(20, 4).Value = "'16.56"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +0.33%  '
$ws.Cells.Item(21, 4).Value = "'493.39"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -1.67%  '
$ws.Cells.Item(22, 4).Value = "'9.18"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -2.11%  '
$ws.Cells.Item(23, 4).Value = "'0.742"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = '  +2.38%  '
$ws.Cells.Item(24, 4).Value = "'0.0000150"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +10.43%  '
$ws.Cells.Item(25, 4).Value = "'85.04"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -1.01%  '
$ws.Cells.Item(26, 4).Value = "'2.39"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -7.15%  '
$ws.Cells.Item(27, 4).Value = "'12.31"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -4.94%  '
$ws.Cells.Item(28, 4).Value = "'10.29"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -6.55%  '
$ws.Cells.Item(29, 5).Value = '  +0.08%  '
$ws.Cells.Item(30, 5).Value = '  +1.22%  '
$ws.Cells.Item(31, 4).Value = "'2.46"
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -2.81%  '
$ws.Cells.Item(32, 4).Value = "'33.29"
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +8.58%  '
$ws.Cells.Item(33, 4).Value = "'7.79"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -4.31%  '
$ws.Cells.Item(34, 5).Value = '  -3.89%  '
$ws.Cells.Item(35, 5).Value = '  -0.23%  '
$ws.Cells.Item(36, 5).Value = '  -3.14%  '
$ws.Cells.Item(37, 5).Value = '  -0.71%  '
$ws.Cells.Item(38, 4).Value = "'5.83"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -5.34%  '
$ws.Cells.Item(39, 4).Value = "'0.332"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -5.58%  '
$ws.Cells.Item(40, 4).Value = "'454.57"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -3.55%  '
$ws.Cells.Item(41, 5).Value = '  -0.96%  '
$ws.Cells.Item(42, 5).Value = '  -2.46%  '
$ws.Cells.Item(43, 4).Value = "'2.91"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -5.24%  '
$ws.Cells.Item(44, 5).Value = '  -2.67%  '
$ws.Cells.Item(45, 4).Value = "'41.48"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -8.48%  '
$ws.Cells.Item(46, 4).Value = '2.851.03'
$ws.Cells.Item(46, 5).Value = '  -3.71%  '
$ws.Cells.Item(47, 4).Value = "'141.39"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +1.59%  '
$ws.Cells.Item(49, 4).Value = "'0.0352"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -2.74%  '
$ws.Cells.Item(50, 4).Value = "'26.30"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -4.09%  '
$ws.Cells.Item(51, 4).Value = "'24.54"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +11.29%  '
